# ------------------------------------------------------------------
# Applies the "Dodano pomysl na trzecia lekcje" edit:
#   1. Splits the "Blockly" word out of the Lekcja-2 bullet item into
#      its own run, wrapped in <w:proofErr> spell-check markers.
#   2. Removes the stray _GoBack bookmark from the "suma dwoch
#      zmiennych" bullet item.
#   3. Appends an empty paragraph, a "Lekcja 3" heading paragraph and
#      a new bulleted paragraph (new numbering list, numId 3) with the
#      "Test - screen z Blockly..." text (incl. proofErr markers and a
#      relocated _GoBack bookmark).
# ------------------------------------------------------------------

$d = $word.ActiveDocument
$app = $word

$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. Split the "Blockly" run in the Lekcja 2 introduction bullet ---------

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Blockly*") {
        $target = $p
        break
    }
}

$targetXml = @"
<w:p $wordNs xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6705FF2E" w14:textId="5AF5C07F" w:rsidR="003D2B96" w:rsidRDefault="003D2B96" w:rsidP="000B1F0D">
  <w:pPr>
    <w:pStyle w:val="Akapitzlist"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Kr&#243;tkie wprowadzenie do zmiennych: utworzenie dw&#243;ch zmiennych, przypisanie im warto&#347;ci i wypisanie (w </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>Blockly</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>)</w:t>
  </w:r>
</w:p>
"@
$target.Range.InsertXML($targetXml)

# --- 2. Remove the stray _GoBack bookmark -----------------------------------

try {
    $gb = $d.Bookmarks("_GoBack")
    $gb.Delete()
} catch {
    # no-op if it is not present
}

# --- 3. Append the new "Lekcja 3" block -------------------------------------

function Append-PlainParagraph([string]$innerXml) {
    $cur = $d.Paragraphs.Last
    $cur.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    $xml = "<w:p $wordNs>" + $innerXml + "</w:p>"
    $newPara.Range.InsertXML($xml)
}

# 3a. a completely empty paragraph
Append-PlainParagraph ""

# 3b. the "Lekcja 3" heading paragraph
Append-PlainParagraph '<w:r><w:t>Lekcja 3</w:t></w:r>'

# 3c. mint a fresh bullet-list numbering definition (numId 3) by applying a
#     bullet list template to a placeholder paragraph, the same way Word
#     itself materializes a brand-new <w:abstractNum>/<w:num> pair the first
#     time a list gallery bullet gets applied to a paragraph.
$cur = $d.Paragraphs.Last
$cur.Range.InsertParagraphAfter()
$placeholder = $d.Paragraphs.Last
$listGallery = $app.ListGalleries.Item(1)
$listTemplate = $listGallery.ListTemplates.Item(1)
$placeholder.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate)

# 3d. replace the placeholder paragraph's content with the final text,
#     preserving the pStyle/numPr that ApplyListTemplateWithLevel produced.
$placeholder = $d.Paragraphs.Last
$bulletXml = @"
<w:p $wordNs>
  <w:pPr>
    <w:pStyle w:val="Akapitzlist"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="3"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Test </w:t>
  </w:r>
  <w:r>
    <w:t>&#8211;</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>screen</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> z </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>Blockly</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> i pytanie: co wypisze program</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
"@
$placeholder.Range.InsertXML($bulletXml)

Write-Output "edit complete"
